$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New resale-number record for 2025-01-21 08:58 (appended as row 35)
$row = 35

# Columns A (date) and D (week number) look numeric/date-like to Excel's
# auto-detection, so force them to be stored as plain text first, then
# drop back to the default "Normal" style so no formatting is left behind.
$textCols = @(1, 4)
foreach ($col in $textCols) {
    $ws.Cells.Item($row, $col).NumberFormat = "@"
}

$ws.Cells.Item($row, 1).Value  = "2025-01-21"
$ws.Cells.Item($row, 2).Value  = "08:58:42"
$ws.Cells.Item($row, 3).Value  = "Tuesday"
$ws.Cells.Item($row, 4).Value  = "03"
$ws.Cells.Item($row, 5).Value  = 126358
$ws.Cells.Item($row, 6).Value  = 142121
$ws.Cells.Item($row, 7).Value  = 168628
$ws.Cells.Item($row, 8).Value  = 158443
$ws.Cells.Item($row, 9).Value  = -1
$ws.Cells.Item($row, 10).Value = 142928
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192372
$ws.Cells.Item($row, 14).Value = 115728
$ws.Cells.Item($row, 15).Value = 45537
$ws.Cells.Item($row, 16).Value = 28487
$ws.Cells.Item($row, 17).Value = 65689
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 48361
$ws.Cells.Item($row, 20).Value = -1

foreach ($col in $textCols) {
    $ws.Cells.Item($row, $col).Style = "Normal"
}
